$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.428.10'
$ws.Range('E2').Value = '  +1.05%  '
$ws.Range('D3').Value = '2.965.49'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.42%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '377.80'
$ws.Range('E5').Value = '  +1.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.53'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.540'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('E8').Value = '  -2.98%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.590'
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.09'
$ws.Range('E10').Value = '  +0.72%  '
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0841'
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').Value = '3.440.10'
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.33'
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.55'
$ws.Range('E15').Value = '  +2.16%  '
$ws.Range('D16').Value = '2.971.23'
$ws.Range('E16').Value = '  +2.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.961'
$ws.Range('E17').Value = '  +2.13%  '
$ws.Range('D18').Value = '51.391.96'
$ws.Range('E18').Value = '  +0.96%  '
$ws.Range('E19').Value = '  +1.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.38'
$ws.Range('E20').Value = '  +1.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.90'
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('D22').Value = '0.0₃0960'
$ws.Range('E22').Value = '  +1.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.28'
$ws.Range('E23').Value = '  +1.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '261.28'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.83'
$ws.Range('E25').Value = '  +5.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.11'
$ws.Range('E26').Value = '  +15.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.53'
$ws.Range('E27').Value = '  +20.34%  '
$ws.Range('E28').Value = '  -1.60%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.113'
$ws.Range('E29').Value = '  +9.07%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.82'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.87'
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.90'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('E34').Value = '  -2.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '50.90'
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  +5.21%  '
$ws.Range('E37').Value = '  +0.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.02'
$ws.Range('E38').Value = '  -0.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.09'
$ws.Range('E39').Value = '  -0.25%  '
$ws.Range('E40').Value = '  -3.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.84'
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.115'
$ws.Range('E42').Value = '  +2.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '124.99'
$ws.Range('E43').Value = '  +4.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.78'
$ws.Range('E44').Value = '  -1.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.289'
$ws.Range('E45').Value = '  +16.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.06'
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('E47').Value = '  +2.69%  '
$ws.Range('D48').Value = '2.028.64'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.21'
$ws.Range('E49').Value = '  +0.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0340'
$ws.Range('E50').Value = '  +9.38%  '
$ws.Range('E51').Value = '  -1.25%  '
